$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value2 = "Datos actualizados a 28 de Marzo de 2020 a las 19:59"
$ws.Cells.Item(4, 2).Value2 = 116448
$ws.Cells.Item(4, 3).Value2 = 12322
$ws.Cells.Item(4, 5).Value2 = 111281
$ws.Cells.Item(22, 1).Value2 = "Israel"
$ws.Cells.Item(22, 2).Value2 = 3619
$ws.Cells.Item(22, 3).Value2 = 584
$ws.Cells.Item(22, 4).Value2 = 89
$ws.Cells.Item(22, 5).Value2 = 3518
$ws.Cells.Item(22, 6).Value2 = 50
$ws.Cells.Item(22, 7).Value2 = 0
$ws.Cells.Item(22, 8).Value2 = 12
$ws.Cells.Item(23, 1).Value2 = "Brasil"
$ws.Cells.Item(23, 2).Value2 = 3477
$ws.Cells.Item(23, 3).Value2 = 60
$ws.Cells.Item(23, 4).Value2 = 6
$ws.Cells.Item(23, 5).Value2 = 3378
$ws.Cells.Item(23, 6).Value2 = 296
$ws.Cells.Item(23, 7).Value2 = 1
$ws.Cells.Item(23, 8).Value2 = 93
$ws.Cells.Item(39, 2).Value2 = 1187
$ws.Cells.Item(39, 3).Value2 = 17
$ws.Cells.Item(39, 5).Value2 = 1154
$ws.Cells.Item(66, 4).Value2 = 30
$ws.Cells.Item(66, 5).Value2 = 374
$ws.Cells.Item(66, 6).Value2 = 4
$ws.Cells.Item(86, 2).Value2 = 224
$ws.Cells.Item(86, 3).Value2 = 1
$ws.Cells.Item(86, 5).Value2 = 196
$ws.Cells.Item(86, 6).Value2 = 16
$ws.Cells.Item(89, 1).Value2 = "Reunion"
$ws.Cells.Item(89, 2).Value2 = 183
$ws.Cells.Item(89, 3).Value2 = 38
$ws.Cells.Item(89, 4).Value2 = 1
$ws.Cells.Item(89, 6).Value2 = 0
$ws.Cells.Item(89, 7).Value2 = 0
$ws.Cells.Item(89, 8).Value2 = 0
$ws.Cells.Item(90, 1).Value2 = "Azerbaiyan"
$ws.Cells.Item(90, 2).Value2 = 182
$ws.Cells.Item(90, 5).Value2 = 163
$ws.Cells.Item(90, 6).Value2 = 23
$ws.Cells.Item(90, 7).Value2 = 1
$ws.Cells.Item(90, 8).Value2 = 4
$ws.Cells.Item(91, 1).Value2 = "Republica de Chipre"
$ws.Cells.Item(91, 2).Value2 = 179
$ws.Cells.Item(91, 3).Value2 = 17
$ws.Cells.Item(91, 4).Value2 = 15
$ws.Cells.Item(91, 5).Value2 = 159
$ws.Cells.Item(91, 8).Value2 = 5
$ws.Cells.Item(92, 1).Value2 = "Vietnam"
$ws.Cells.Item(92, 2).Value2 = 174
$ws.Cells.Item(92, 4).Value2 = 21
$ws.Cells.Item(92, 5).Value2 = 153
$ws.Cells.Item(92, 6).Value2 = 3
$ws.Cells.Item(93, 1).Value2 = "Islas Feroe"
$ws.Cells.Item(93, 2).Value2 = 155
$ws.Cells.Item(93, 3).Value2 = 11
$ws.Cells.Item(93, 4).Value2 = 54
$ws.Cells.Item(93, 5).Value2 = 101
$ws.Cells.Item(93, 6).Value2 = 2
$ws.Cells.Item(94, 1).Value2 = "Oman"
$ws.Cells.Item(94, 2).Value2 = 152
$ws.Cells.Item(94, 3).Value2 = 21
$ws.Cells.Item(94, 4).Value2 = 23
$ws.Cells.Item(94, 5).Value2 = 129
$ws.Cells.Item(94, 6).Value2 = 0
$ws.Cells.Item(95, 1).Value2 = "Malta"
$ws.Cells.Item(95, 2).Value2 = 149
$ws.Cells.Item(95, 3).Value2 = 10
$ws.Cells.Item(95, 4).Value2 = 2
$ws.Cells.Item(95, 5).Value2 = 147
$ws.Cells.Item(95, 6).Value2 = 1
$ws.Cells.Item(105, 1).Value2 = "Guadalupe"
$ws.Cells.Item(105, 2).Value2 = 102
$ws.Cells.Item(105, 3).Value2 = 29
$ws.Cells.Item(105, 4).Value2 = 17
$ws.Cells.Item(105, 5).Value2 = 83
$ws.Cells.Item(105, 6).Value2 = 4
$ws.Cells.Item(105, 7).Value2 = 1
$ws.Cells.Item(105, 8).Value2 = 2
$ws.Cells.Item(106, 1).Value2 = "Costa de Marfil"
$ws.Cells.Item(106, 2).Value2 = 101
$ws.Cells.Item(106, 4).Value2 = 3
$ws.Cells.Item(106, 5).Value2 = 98
$ws.Cells.Item(106, 6).Value2 = 0
$ws.Cells.Item(107, 1).Value2 = "Camboya"
$ws.Cells.Item(107, 2).Value2 = 99
$ws.Cells.Item(107, 3).Value2 = 0
$ws.Cells.Item(107, 4).Value2 = 13
$ws.Cells.Item(107, 5).Value2 = 86
$ws.Cells.Item(107, 6).Value2 = 1
$ws.Cells.Item(107, 8).Value2 = 0
$ws.Cells.Item(108, 1).Value2 = "Estado de Palestina"
$ws.Cells.Item(108, 2).Value2 = 97
$ws.Cells.Item(108, 3).Value2 = 6
$ws.Cells.Item(108, 4).Value2 = 18
$ws.Cells.Item(108, 5).Value2 = 78
$ws.Cells.Item(108, 6).Value2 = 0
$ws.Cells.Item(108, 7).Value2 = 0
$ws.Cells.Item(108, 8).Value2 = 1
$ws.Cells.Item(119, 1).Value2 = "Ruanda"
$ws.Cells.Item(119, 2).Value2 = 60
$ws.Cells.Item(119, 3).Value2 = 6
$ws.Cells.Item(119, 5).Value2 = 60
$ws.Cells.Item(120, 1).Value2 = "Kirguistan"
$ws.Cells.Item(120, 3).Value2 = 0
$ws.Cells.Item(120, 4).Value2 = 0
$ws.Cells.Item(120, 5).Value2 = 58
$ws.Cells.Item(120, 7).Value2 = 0
$ws.Cells.Item(120, 8).Value2 = 0
$ws.Cells.Item(121, 1).Value2 = "Consejo Danes para los Refugiados"
$ws.Cells.Item(121, 2).Value2 = 58
$ws.Cells.Item(121, 3).Value2 = 7
$ws.Cells.Item(121, 4).Value2 = 2
$ws.Cells.Item(121, 5).Value2 = 50
$ws.Cells.Item(121, 7).Value2 = 3
$ws.Cells.Item(121, 8).Value2 = 6
$ws.Cells.Item(122, 1).Value2 = "Liechtenstein"
$ws.Cells.Item(122, 3).Value2 = 0
$ws.Cells.Item(122, 4).Value2 = 0
$ws.Cells.Item(122, 5).Value2 = 56
$ws.Cells.Item(122, 6).Value2 = 0
$ws.Cells.Item(122, 8).Value2 = 0
$ws.Cells.Item(123, 1).Value2 = "Paraguay"
$ws.Cells.Item(123, 3).Value2 = 4
$ws.Cells.Item(123, 4).Value2 = 1
$ws.Cells.Item(123, 5).Value2 = 52
$ws.Cells.Item(123, 6).Value2 = 1
$ws.Cells.Item(123, 8).Value2 = 3
$ws.Cells.Item(124, 1).Value2 = "Gibraltar"
$ws.Cells.Item(124, 2).Value2 = 56
$ws.Cells.Item(124, 3).Value2 = 1
$ws.Cells.Item(124, 4).Value2 = 14
$ws.Cells.Item(124, 5).Value2 = 42
$ws.Cells.Item(153, 1).Value2 = "Dominica"
$ws.Cells.Item(154, 1).Value2 = "San Martin (Parte Francesa)"
$ws.Cells.Item(155, 1).Value2 = "Bahamas"
$ws.Cells.Item(155, 4).Value2 = 1
$ws.Cells.Item(155, 8).Value2 = 0
$ws.Cells.Item(156, 1).Value2 = "Niger"
$ws.Cells.Item(156, 4).Value2 = 0
$ws.Cells.Item(156, 8).Value2 = 1
$ws.Cells.Item(159, 1).Value2 = "Surinam"
$ws.Cells.Item(160, 1).Value2 = "Haiti"
$ws.Cells.Item(161, 1).Value2 = "Birmania"
$ws.Cells.Item(162, 1).Value2 = "Laos"
$ws.Cells.Item(162, 3).Value2 = 2
$ws.Cells.Item(164, 1).Value2 = "Mozambique"
$ws.Cells.Item(164, 3).Value2 = 1
$ws.Cells.Item(165, 1).Value2 = "Guyana"
$ws.Cells.Item(165, 3).Value2 = 3
$ws.Cells.Item(166, 1).Value2 = "Islas Caimanes"
$ws.Cells.Item(166, 4).Value2 = 0
$ws.Cells.Item(166, 5).Value2 = 7
$ws.Cells.Item(166, 8).Value2 = 1
$ws.Cells.Item(167, 1).Value2 = "Namibia"
$ws.Cells.Item(167, 5).Value2 = 6
$ws.Cells.Item(167, 8).Value2 = 0
$ws.Cells.Item(168, 1).Value2 = "Curazao"
$ws.Cells.Item(168, 2).Value2 = 8
$ws.Cells.Item(168, 4).Value2 = 2
$ws.Cells.Item(168, 5).Value2 = 5
$ws.Cells.Item(168, 8).Value2 = 1
$ws.Cells.Item(169, 1).Value2 = "Antigua y Barbuda"
$ws.Cells.Item(171, 1).Value2 = "Granada"
$ws.Cells.Item(171, 3).Value2 = 0
$ws.Cells.Item(171, 5).Value2 = 7
$ws.Cells.Item(171, 8).Value2 = 0
$ws.Cells.Item(172, 1).Value2 = "Zimbabue"
$ws.Cells.Item(172, 3).Value2 = 2
$ws.Cells.Item(173, 1).Value2 = "Gabon"
$ws.Cells.Item(173, 2).Value2 = 7
$ws.Cells.Item(173, 8).Value2 = 1
$ws.Cells.Item(175, 1).Value2 = "Eritrea"
$ws.Cells.Item(176, 1).Value2 = "Benin"
$ws.Cells.Item(176, 2).Value2 = 6
$ws.Cells.Item(176, 5).Value2 = 6
$ws.Cells.Item(177, 1).Value2 = "San Bartolome"
$ws.Cells.Item(178, 1).Value2 = "Fiyi"
$ws.Cells.Item(179, 1).Value2 = "Montserrat"
$ws.Cells.Item(180, 1).Value2 = "Siria"
$ws.Cells.Item(180, 3).Value2 = 0
$ws.Cells.Item(181, 1).Value2 = "Mauritania"
$ws.Cells.Item(181, 3).Value2 = 2
$ws.Cells.Item(181, 4).Value2 = 0
$ws.Cells.Item(181, 5).Value2 = 5
$ws.Cells.Item(184, 1).Value2 = "Nepal"
$ws.Cells.Item(184, 3).Value2 = 1
$ws.Cells.Item(184, 4).Value2 = 1
$ws.Cells.Item(184, 8).Value2 = 0
